# Applies a cyclic rotation of the weekly date groups for rows 2-7
# (Fruta, Vega Monumental Concepcion - Mora) in the active worksheet.
#
# Rows 2-3 (Primera/Segunda) <- values formerly in rows 6-7
# Rows 4-5 (Primera/Segunda) <- same values, only the date moves 44559 -> 44574
# Rows 6-7 (Primera/Segunda) <- values formerly in rows 2-3 (== rows 4-5), with date 44559

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Primera)
$ws.Range("D2").Value = 44223
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 3500
$ws.Range("O2").Value = 4000
$ws.Range("P2").Value = 3750
$ws.Range("S2").Value = 1875

# Row 3 (Segunda)
$ws.Range("D3").Value = 44223
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("S3").Value = 1500

# Row 4 (Primera) - only date changes
$ws.Range("D4").Value = 44574

# Row 5 (Segunda) - only date changes
$ws.Range("D5").Value = 44574

# Row 6 (Primera)
$ws.Range("D6").Value = 44559
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 6000
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 6500
$ws.Range("S6").Value = 3250

# Row 7 (Segunda)
$ws.Range("D7").Value = 44559
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 5000
$ws.Range("P7").Value = 5000
$ws.Range("S7").Value = 2500
